$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header info strings
$ws.Range("C4").Value = "Route : O2"
$ws.Range("E4").Value = "From : 01-01-2019"
$ws.Range("G4").Value = "To : 03-05-2019"

# Insert a new row before the old "Print taken by" row (currently row 9)
$ws.Rows.Item(9).Insert()

# Fill in the new data row (row 8)
$ws.Cells.Item(8, 1).Value = 1
$ws.Cells.Item(8, 2).Value = "O2-2-1"
$ws.Cells.Item(8, 3).Value = "2:30 PM"
$ws.Cells.Item(8, 4).Value = 1
$ws.Cells.Item(8, 5).Value = 8
$ws.Cells.Item(8, 6).Value = 13305
$ws.Cells.Item(8, 7).Value = 1234
$ws.Cells.Item(8, 8).Value = 0
$ws.Cells.Item(8, 9).Value = -1234
$ws.Cells.Item(8, 10).Value = 221
$ws.Cells.Item(8, 11).Value = 956

# Update the "Print taken at" text on the shifted footer row (row 10)
$ws.Range("F10").Value = "Print taken at : 03-05-2019 16:35:12"

$ws.Range("F10").Select()
